$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTAS POR GRUPO")

# --- Step 1: Insert a new column at F (6), shifting GRIFERIAS..SAL SOLUBLE (old F:N) to G:O ---
$ws.Columns.Item(6).Insert(-4161)

# Set header and data for the newly inserted column F ("GRANITO").
# The Insert() operation already copied the formatting (styles) from the row across,
# so we only need to populate the values.
$ws.Range("F1").Value = "GRANITO"
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("F7").Value = "0 de 5"

# Set the column width for the new column F to match the target layout.
$ws.Columns.Item(6).ColumnWidth = 12.166666666666666

# --- Step 2: Append three new trailing columns P, Q, R ---
# Copy formatting from the last existing column (O, formerly N "SAL SOLUBLE") into each
# new column/row so the cell styles (header bold/border, currency number format,
# centered "x de 5" text) match the rest of the table, then set the actual values.

$headers = @("NO RESURTIBLES", "PANELES PVC", "PANELES PU")
$cols = @("P", "Q", "R")
$colNums = @(16, 17, 18)
$widths = @(19.166666666666668, 16.166666666666668, 15.166666666666666)

for ($i = 0; $i -lt 3; $i++) {
    $col = $cols[$i]

    $ws.Range("O1").Copy($ws.Range($col + "1"))
    $ws.Range($col + "1").Value = $headers[$i]

    for ($r = 2; $r -le 6; $r++) {
        $ws.Range("O" + $r).Copy($ws.Range($col + $r))
        $ws.Range($col + $r).Value = 0
    }

    $ws.Range("O7").Copy($ws.Range($col + "7"))
    $ws.Range($col + "7").Value = "0 de 5"

    $ws.Columns.Item($colNums[$i]).ColumnWidth = $widths[$i]
}
